$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; this shifts the existing rows 17-21
# down to 18-22 (with all their original values/formatting intact).
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new data record.
$ws.Cells.Item(17, 1).Value = 12
$ws.Cells.Item(17, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44510
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = 100112028
$ws.Cells.Item(17, 7).Value = "Sandia"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 250
$ws.Cells.Item(17, 11).Value = 800
$ws.Cells.Item(17, 12).Value = 800
$ws.Cells.Item(17, 13).Value = 800
$ws.Cells.Item(17, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(17, 15).Value = "Perú"
$ws.Cells.Item(17, 16).Value = 800
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = "Hortaliza"
